# "Generate Report for Handoff"
# Updates the localization-status report: the Overview/zh-cn/de-de sheets
# move from "In Translation" to "Ready for handoff", the associated
# timestamps advance, and the now-wider "Status" columns are resized to fit.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Timestamps advance to reflect the new handoff generation --------------
$overview.Range("G2").Value = "2016-08-16 06:34:39"
$dede.Range("H2").Value     = "2016-08-16 06:34:39"
$zhcn.Range("H2").Value     = "2016-08-16 06:34:34"

# --- Widen the Status columns so the new, longer text fits -----------------
$overview.Columns.Item(5).ColumnWidth = 16.38265482584637
$overview.Columns.Item(6).ColumnWidth = 16.38265482584637
$zhcn.Columns.Item(3).ColumnWidth     = 16.38265482584637
$dede.Columns.Item(3).ColumnWidth     = 16.38265482584637
